# Plant-Data-Tidy.xlsx: normalize the "Trial" labels on Survival_cdf.Data
# (Sheet1) by removing the internal space after "Test" (and, for the
# 3rd/4th trials, the trailing space as well), e.g.
#   "Test 1LL " -> "Test1LL "
#   "Test 2LL " -> "Test2LL "
#   "Test 3LL " -> "Test3LL"
#   "Test 4LL " -> "Test4LL"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Survival_cdf.Data")

$ws.Range("A2:A9").Value   = "Test1LL "
$ws.Range("A10:A16").Value = "Test2LL "
$ws.Range("A17:A23").Value = "Test3LL"
$ws.Range("A24:A29").Value = "Test4LL"
